$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("D2")
$cell.NumberFormat = "@"
$cell.Value = "26.247.48"
$cell.Style = "Normal"
$ws.Range("E2").Value = "  +0.48%  "

$cell = $ws.Range("D3")
$cell.NumberFormat = "@"
$cell.Value = "1.655.11"
$cell.Style = "Normal"
$ws.Range("E3").Value = "  -0.08%  "

$cell = $ws.Range("D4")
$cell.NumberFormat = "@"
$cell.Value = "1.005"
$cell.Style = "Normal"
$ws.Range("E4").Value = "  +0.44%  "

$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "219.11"
$cell.Style = "Normal"

$ws.Range("E6").Value = "  -0.42%  "

$cell = $ws.Range("D7")
$cell.NumberFormat = "@"
$cell.Value = "1.005"
$cell.Style = "Normal"
$ws.Range("E7").Value = "  +0.43%  "

$cell = $ws.Range("D8")
$cell.NumberFormat = "@"
$cell.Value = "0.2651"
$cell.Style = "Normal"
$ws.Range("E8").Value = "  +1.09%  "

$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = "0.06315"
$cell.Style = "Normal"

$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = "21.39"
$cell.Style = "Normal"
$ws.Range("E10").Value = "  +2.69%  "

$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = "0.07739"
$cell.Style = "Normal"
$ws.Range("E11").Value = "  -0.18%  "

$cell = $ws.Range("D12")
$cell.NumberFormat = "@"
$cell.Value = "1.676.54"
$cell.Style = "Normal"
$ws.Range("E12").Value = "  +1.14%  "

$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = "4.427"
$cell.Style = "Normal"
$ws.Range("E13").Value = "  -0.52%  "

$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = "0.5454"
$cell.Style = "Normal"
$ws.Range("E14").Value = "  -1.05%  "

$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = "0.0₅8222"
$cell.Style = "Normal"
$ws.Range("E15").Value = "  -0.60%  "

$ws.Range("E16").Value = "  -0.51%  "

$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = "26.246.06"
$cell.Style = "Normal"
$ws.Range("E17").Value = "  +0.45%  "

$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = "1.005"
$cell.Style = "Normal"
$ws.Range("E18").Value = "  +0.42%  "

$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = "4.688"
$cell.Style = "Normal"
$ws.Range("E19").Value = "  -1.49%  "

$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = "191.35"
$cell.Style = "Normal"
$ws.Range("E20").Value = "  +0.54%  "

$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = "10.18"
$cell.Style = "Normal"
$ws.Range("E21").Value = "  -1.20%  "

$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = "6.200"
$cell.Style = "Normal"
$ws.Range("E22").Value = "  -2.44%  "

$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = "1.007"
$cell.Style = "Normal"
$ws.Range("E23").Value = "  +0.52%  "

$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = "139.06"
$cell.Style = "Normal"
$ws.Range("E24").Value = "  -2.81%  "

$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = "0.1246"
$cell.Style = "Normal"

$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value = "7.294"
$cell.Style = "Normal"
$ws.Range("E26").Value = "  -1.47%  "

$ws.Range("E27").Value = "  +0.40%  "

$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = "1.418"
$cell.Style = "Normal"
$ws.Range("E28").Value = "  -1.00%  "

$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = "0.06059"
$cell.Style = "Normal"
$ws.Range("E29").Value = "  +0.29%  "

$ws.Range("E30").Value = "  +1.83%  "

$cell = $ws.Range("D31")
$cell.NumberFormat = "@"
$cell.Value = "3.551"
$cell.Style = "Normal"
$ws.Range("E31").Value = "  +1.37%  "

$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = "3.361"
$cell.Style = "Normal"
$ws.Range("E32").Value = "  -1.60%  "

$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = "1.655"
$cell.Style = "Normal"
$ws.Range("E33").Value = "  -0.14%  "

$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = "0.9870"
$cell.Style = "Normal"

$cell = $ws.Range("D35")
$cell.NumberFormat = "@"
$cell.Value = "2.418"
$cell.Style = "Normal"
$ws.Range("E35").Value = "  +0.78%  "

$cell = $ws.Range("D36")
$cell.NumberFormat = "@"
$cell.Value = "2.774"
$cell.Style = "Normal"
$ws.Range("E36").Value = "  +0.63%  "

$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = "0.5958"
$cell.Style = "Normal"
$ws.Range("E37").Value = "  +5.49%  "

$cell = $ws.Range("D38")
$cell.NumberFormat = "@"
$cell.Value = "0.01597"
$cell.Style = "Normal"
$ws.Range("E38").Value = "  -0.42%  "

$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value = "5.966"
$cell.Style = "Normal"
$ws.Range("E39").Value = "  +1.10%  "

$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = "1.065.03"
$cell.Style = "Normal"
$ws.Range("E40").Value = "  +3.49%  "

$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = "0.8600"
$cell.Style = "Normal"
$ws.Range("E41").Value = "  +0.71%  "

$ws.Range("E42").Value = "  +0.31%  "

$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = "99.94"
$cell.Style = "Normal"

$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = "1.798.25"
$cell.Style = "Normal"
$ws.Range("E44").Value = "  -0.31%  "

$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = "0.0₈109"
$cell.Style = "Normal"
$ws.Range("E45").Value = "  +2.94%  "

$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = "57.35"
$cell.Style = "Normal"
$ws.Range("E46").Value = "  +2.49%  "

$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = "1.004"
$cell.Style = "Normal"
$ws.Range("E47").Value = "  +0.32%  "

$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = "8.069"
$cell.Style = "Normal"
$ws.Range("E48").Value = "  -0.44%  "

$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = "0.05180"
$cell.Style = "Normal"
$ws.Range("E49").Value = "  +0.46%  "

$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = "1.471"
$cell.Style = "Normal"
$ws.Range("E50").Value = "  +5.40%  "

$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = "0.4232"
$cell.Style = "Normal"
$ws.Range("E51").Value = "  +0.40%  "
